$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2 through 13
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25)
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
